# "new changes in ops (ordercreation & orderpage & order form)"
# - Rename "Emp ID-Order Assigned" -> "Typist" (col C header)
# - Rename "Assignee_QA" -> "Typist QC" (col D header)
# - Update sample Typist id  SIPL5316 -> SIPL0102
# - Update sample Typist QC id SIPL5688 -> SIPL5317
# - Update Status sample value WIP -> Typing
# - Column C's data cell loses its top border (no longer fully boxed on top)
# - Minor cosmetic: active selection moves to E5, column widths re-autosized

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "Typist"
$ws.Range("D1").Value = "Typist QC"

# --- Data row (row 2) ---
$ws.Range("C2").Value = "SIPL0102"
$ws.Range("D2").Value = "SIPL5317"
$ws.Range("L2").Value = "Typing"

# --- Border tweak: C2 keeps its left/right/bottom thin border but drops the
# top edge (it no longer needs to be boxed against the header row). ---
$ws.Range("C2").Borders.Item(8).LineStyle = 0

# --- Column widths re-fit around the new, shorter header/values ---
$ws.Columns.Item(3).ColumnWidth = 7.42
$ws.Columns.Item(8).ColumnWidth = 16.6

# --- Selection moved before save ---
$ws.Range("E5").Select()
